# The reviews_count column (E) is no longer produced by the scraper, so
# remove it entirely and let the remaining columns (reviews_average,
# latitude, longitude, is_permanently_closed, gmaps_link,
# latest_review_date) shift one position to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("E").Delete()
